$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1356.6666
$ws.Range("I98").Value = 1356.6666
$ws.Range("K98").Value = 1356.6666
$ws.Range("M98").Value = 141.3334
$ws.Range("H100").Value = 14929341
$ws.Range("I100").Value = 20834684
$ws.Range("J100").Value = 1431416.1
$ws.Range("K100").Value = 20834684
$ws.Range("L100").Value = 1431416.1
$ws.Range("M100").Value = -20834143
$ws.Range("N100").Value = -1432498.1
$ws.Range("H103").Value = 801
$ws.Range("I103").Value = 1004
$ws.Range("J103").Value = 598
$ws.Range("K103").Value = 3012
$ws.Range("L103").Value = 1794
$ws.Range("M103").Value = -2426
$ws.Range("N103").Value = -2966
$ws.Range("H107").Value = 270.57144
$ws.Range("I107").Value = 267.42856
$ws.Range("K107").Value = 267.42856
$ws.Range("M107").Value = 1652.57144
$ws.Range("H122").Value = 1356.6666
$ws.Range("I122").Value = 1356.6666
$ws.Range("K122").Value = 4069.9998
$ws.Range("M122").Value = -1619.9998
$ws.Range("H129").Value = 1115.1786
$ws.Range("I129").Value = 476.125
$ws.Range("K129").Value = 1428.375
$ws.Range("M129").Value = 3571.625
$ws.Range("H132").Value = 58670.668
$ws.Range("I132").Value = 65748.31
$ws.Range("J132").Value = 2049.5
$ws.Range("K132").Value = 197244.93
$ws.Range("L132").Value = 6148.5
$ws.Range("M132").Value = -194714.93
$ws.Range("N132").Value = -11208.5
$ws.Range("H137").Value = 1191.7413
$ws.Range("I137").Value = 1192.8857
$ws.Range("J137").Value = 1190
$ws.Range("K137").Value = 3578.6571
$ws.Range("L137").Value = 3570
$ws.Range("M137").Value = -1028.6571
$ws.Range("N137").Value = -8670

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3246.5762
$ws.Range("I32").Value = 2090.4119
$ws.Range("K32").Value = 2090.4119
$ws.Range("M32").Value = -1803.4119
$ws.Range("H61").Value = 3988.2153
$ws.Range("I61").Value = 4816.636
$ws.Range("K61").Value = 4816.636
$ws.Range("M61").Value = -4604.636
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H97").Value = 1564.9
$ws.Range("I97").Value = 1331.125
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 1331.125
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -835.125
$ws.Range("N97").Value = -3492
$ws.Range("H101").Value = 62500
$ws.Range("J101").Value = 62500
$ws.Range("L101").Value = 62500
$ws.Range("N101").Value = -68990
$ws.Range("H102").Value = 2200
$ws.Range("I102").Value = 2200
$ws.Range("K102").Value = 2200
$ws.Range("M102").Value = -578
$ws.Range("H104").Value = 32593.75
$ws.Range("J104").Value = 32593.75
$ws.Range("L104").Value = 32593.75
$ws.Range("N104").Value = -39581.75
$ws.Range("H122").Value = 1690.2727
$ws.Range("I122").Value = 1671.6052
$ws.Range("J122").Value = 1732
$ws.Range("K122").Value = 5014.8156
$ws.Range("L122").Value = 5196
$ws.Range("M122").Value = -2564.8156
$ws.Range("N122").Value = -10096
$ws.Range("H136").Value = 3988.2153
$ws.Range("I136").Value = 4816.636
$ws.Range("K136").Value = 14449.908
$ws.Range("M136").Value = -11899.908

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 594.875
$ws.Range("I94").Value = 580.5333000000001
$ws.Range("J94").Value = 810
$ws.Range("K94").Value = 580.5333000000001
$ws.Range("L94").Value = 810
$ws.Range("M94").Value = -129.5333000000001
$ws.Range("N94").Value = -1712
$ws.Range("H99").Value = 1476.8422
$ws.Range("I99").Value = 895.55554
$ws.Range("K99").Value = 895.55554
$ws.Range("M99").Value = 602.44446
$ws.Range("H105").Value = 1660
$ws.Range("I105").Value = 1603.3334
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1603.3334
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 143.6666
$ws.Range("N105").Value = -5494

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1433.85
$ws.Range("I58").Value = 1197.1538
$ws.Range("K58").Value = 1197.1538
$ws.Range("M58").Value = -994.1538
$ws.Range("H134").Value = 4686.276
$ws.Range("I134").Value = 5391.75
$ws.Range("K134").Value = 16175.25
$ws.Range("M134").Value = -13640.25
$ws.Range("H136").Value = 1433.85
$ws.Range("I136").Value = 1197.1538
$ws.Range("K136").Value = 3591.4614
$ws.Range("M136").Value = -1041.4614

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 214.625
$ws.Range("I2").Value = 280.875
$ws.Range("J2").Value = 148.375
$ws.Range("K2").Value = 1685.25
$ws.Range("L2").Value = 890.25
$ws.Range("M2").Value = -1572.25
$ws.Range("N2").Value = -1116.25
$ws.Range("H12").Value = 43478348
$ws.Range("I12").Value = 200000000
$ws.Range("J12").Value = 112.111115
$ws.Range("K12").Value = 600000000
$ws.Range("L12").Value = 336.333345
$ws.Range("M12").Value = -599999827
$ws.Range("N12").Value = -682.333345
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 3318
$ws.Range("J34").Value = 3318
$ws.Range("L34").Value = 9954
$ws.Range("N34").Value = -10122
$ws.Range("H39").Value = 17512.875
$ws.Range("J39").Value = 5728.5713
$ws.Range("L39").Value = 17185.7139
$ws.Range("N39").Value = -17773.7139
$ws.Range("H55").Value = 2993.3333
$ws.Range("J55").Value = 3000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354
$ws.Range("H68").Value = 991.5567
$ws.Range("I68").Value = 704.50946
$ws.Range("J68").Value = 1337.3182
$ws.Range("K68").Value = 2113.52838
$ws.Range("L68").Value = 4011.9546
$ws.Range("M68").Value = -1302.52838
$ws.Range("N68").Value = -5633.9546
$ws.Range("H71").Value = 991.5567
$ws.Range("I71").Value = 704.50946
$ws.Range("J71").Value = 1337.3182
$ws.Range("K71").Value = 6340.58514
$ws.Range("L71").Value = 12035.8638
$ws.Range("M71").Value = -2284.58514
$ws.Range("N71").Value = -20147.8638

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1987.45
$ws.Range("I97").Value = 1788.5
$ws.Range("J97").Value = 2451.6667
$ws.Range("K97").Value = 1788.5
$ws.Range("L97").Value = 2451.6667
$ws.Range("M97").Value = -1292.5
$ws.Range("N97").Value = -3443.6667
$ws.Range("H122").Value = 62222.51
$ws.Range("I122").Value = 77359
$ws.Range("J122").Value = 12272.1
$ws.Range("K122").Value = 232077
$ws.Range("L122").Value = 36816.3
$ws.Range("M122").Value = -229627
$ws.Range("N122").Value = -41716.3
$ws.Range("H126").Value = 2194.5518
$ws.Range("I126").Value = 1860.3636
$ws.Range("J126").Value = 3244.8572
$ws.Range("K126").Value = 5581.0908
$ws.Range("L126").Value = 9734.571599999999
$ws.Range("M126").Value = -3111.0908
$ws.Range("N126").Value = -14674.5716

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 22500
$ws.Range("J98").Value = 22500
$ws.Range("L98").Value = 22500
$ws.Range("N98").Value = -28490
$ws.Range("H100").Value = 1851
$ws.Range("I100").Value = 1851
$ws.Range("K100").Value = 1851
$ws.Range("M100").Value = -1310
$ws.Range("I104").Value = 20000
$ws.Range("K104").Value = 20000
$ws.Range("M104").Value = -16506

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 7943.75
$ws.Range("J45").Value = 7943.75
$ws.Range("L45").Value = 7943.75
$ws.Range("N45").Value = -8925.75
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2341.3333
$ws.Range("I132").Value = 1582.28
$ws.Range("J132").Value = 3166.3914
$ws.Range("K132").Value = 4746.84
$ws.Range("L132").Value = 9499.174199999999
$ws.Range("M132").Value = -2216.84
$ws.Range("N132").Value = -14559.1742
$ws.Range("H136").Value = 4632357
$ws.Range("I136").Value = 12346165
$ws.Range("J136").Value = 4072.2
$ws.Range("K136").Value = 37038495
$ws.Range("L136").Value = 12216.6
$ws.Range("M136").Value = -37035945
$ws.Range("N136").Value = -17316.6
